$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.700.89"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.846.70"
$ws.Range("E3").Value = "  -1.00%  "
$ws.Range("E4").Value = "  -1.97%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "320.06"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.013"
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4308"
$ws.Range("E7").Value = "  -2.59%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3738"
$ws.Range("E8").Value = "  -1.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07352"
$ws.Range("E9").Value = "  -1.51%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8769"
$ws.Range("E10").Value = "  -1.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.60"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.855.39"
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.717"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.441"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07129"
$ws.Range("E15").Value = "  -1.01%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.12"
$ws.Range("E16").Value = "  +4.47%  "
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008989"
$ws.Range("E18").Value = "  -1.54%  "
$ws.Range("E19").Value = "  -1.91%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.44"
$ws.Range("E20").Value = "  -0.83%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.711.74"
$ws.Range("E21").Value = "  -0.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.223"
$ws.Range("E22").Value = "  -1.62%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.09"
$ws.Range("E23").Value = "  -2.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.086.97"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.012"
$ws.Range("E25").Value = "  -0.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.72"
$ws.Range("E26").Value = "  -2.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.61"
$ws.Range("E27").Value = "  -1.29%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.149"
$ws.Range("E28").Value = "  +7.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.376"
$ws.Range("E29").Value = "  +0.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.87"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08964"
$ws.Range("E31").Value = "  -0.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.229"
$ws.Range("E32").Value = "  -0.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7756"
$ws.Range("E33").Value = "  -0.58%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.552"
$ws.Range("E34").Value = "  -0.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.917"
$ws.Range("E35").Value = "  -3.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.014"
$ws.Range("E36").Value = "  -1.92%  "
$ws.Range("E37").Value = "  -1.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05334"
$ws.Range("E38").Value = "  -0.30%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01972"
$ws.Range("E39").Value = "  -1.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.277"
$ws.Range("E40").Value = "  +5.62%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.882"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5135"
$ws.Range("E42").Value = "  -1.79%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1684"
$ws.Range("E43").Value = "  -0.74%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.774"
$ws.Range("E44").Value = "  +0.89%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.72"
$ws.Range("E45").Value = "  +0.27%  "
$ws.Range("B46").Value = "Quant"
$ws.Range("C46").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "109.55"
$ws.Range("E46").Value = "  -1.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4740"
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06495"
$ws.Range("E48").Value = "  -3.99%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.014"
$ws.Range("E49").Value = "  -1.99%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.694"
$ws.Range("E50").Value = "  -1.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.847"
$ws.Range("E51").Value = "  -4.22%  "
